# Weekly update: a new Berenjena (Vega Monumental Concepción) price record
# is inserted as the new row 68, pushing the previously-existing rows
# 68-94 down to 69-95 (the data itself is unchanged, only its row position
# shifts). The sheet's used range therefore grows from A1:R94 to A1:R95.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 68, shifting rows 68:94 -> 69:95.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new weekly record.
$ws.Cells.Item(68, 1).Value  = 11
$ws.Cells.Item(68, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(68, 3).Value  = "Bíobío"
$ws.Cells.Item(68, 4).Value  = 44726
$ws.Cells.Item(68, 5).Value  = 8
$ws.Cells.Item(68, 6).Value  = 100112001
$ws.Cells.Item(68, 7).Value  = "Berenjena"
$ws.Cells.Item(68, 8).Value  = "Sin especificar"
$ws.Cells.Item(68, 9).Value  = "Primera"
$ws.Cells.Item(68, 10).Value = 160
$ws.Cells.Item(68, 11).Value = 7000
$ws.Cells.Item(68, 12).Value = 7500
$ws.Cells.Item(68, 13).Value = 7250
$ws.Cells.Item(68, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(68, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(68, 16).Value = 121
$ws.Cells.Item(68, 17).Value = 60
$ws.Cells.Item(68, 18).Value = "Hortaliza"
